$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23199474811554
$ws.Range("B1").Value = 2.496006488800049
$ws.Range("C1").Value = 4.196475028991699
$ws.Range("D1").Value = 2.721020936965942
$ws.Range("E1").Value = 1.08258068561554
